$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.05 = 28329.81 pesos`n✅ 28329.81 pesos = 7.01 = 967.22 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the "tasas" sheet rates ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 141.9
$ws2.Range("O10").Value = 4020
$ws2.Range("N12").Value = 4042
$ws2.Range("O12").Value = 138
